# Form_A_Duplicate_1.docx template updates:
#  - {shareholderName_N} placeholders (used for the share-certificate
#    listing blocks) become {shareholderNameCertificate_N}.
#  - the lone {companyOldName2} placeholder is replaced by {companyName}
#    followed by a line break and a conditional {companyOldName} block.

$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-All($findText, $replaceText) {
    $d.Content.Find.Execute(
        $findText,    # FindText
        $false,       # MatchCase
        $false,       # MatchWholeWord
        $false,       # MatchWildcards
        $false,       # MatchSoundsLike
        $false,       # MatchAllWordForms
        $true,        # Forward
        $wdFindContinue, # Wrap
        $false,       # Format
        $replaceText, # ReplaceWith
        $wdReplaceAll # Replace
    ) | Out-Null
}

# {#hasShareholder_1}{shareholderName_1} ...  -> {shareholderNameCertificate_1}
Replace-All "{#hasShareholder_1}{shareholderName_1} " "{#hasShareholder_1}{shareholderNameCertificate_1} "
Replace-All "{#hasShareholder_1}{shareholderName_1};" "{#hasShareholder_1}{shareholderNameCertificate_1};"

# {/hasShareholder_1} {#hasShareholder_2}{shareholderName_2} ... -> {shareholderNameCertificate_2}
Replace-All "{/hasShareholder_1} {#hasShareholder_2}{shareholderName_2} " "{/hasShareholder_1} {#hasShareholder_2}{shareholderNameCertificate_2} "
Replace-All "{/hasShareholder_1} {#hasShareholder_2}{shareholderName_2};" "{/hasShareholder_1} {#hasShareholder_2}{shareholderNameCertificate_2};"

# {/hasShareholder_2} {#hasShareholder_3}{shareholderName_3} ... -> {shareholderNameCertificate_3}
Replace-All "{/hasShareholder_2} {#hasShareholder_3}{shareholderName_3} " "{/hasShareholder_2} {#hasShareholder_3}{shareholderNameCertificate_3} "
Replace-All "{/hasShareholder_2} {#hasShareholder_3}{shareholderName_3};" "{/hasShareholder_2} {#hasShareholder_3}{shareholderNameCertificate_3};"

# {companyOldName2} -> {companyName}<line break>{#hasCompanyOldName}[{companyOldName}]{/hasCompanyOldName}
Replace-All "{companyOldName2}" "{companyName}^l{#hasCompanyOldName}[{companyOldName}]{/hasCompanyOldName}"
